# Edit script for "Niels Mejia Final Homework Presentation.pptx"
# Implements the diff for Slide 1 (sldId 256 / cId 3338301109):
#   - NATURAL GAS / FUEL GAS / FLEXIGAS / LIGHT NAPHTHA boxes get updated
#     labels (with numeric values), smaller font (12pt), and repositioned /
#     resized text boxes.
#   - The connector between FUEL GAS and FLEXIGAS moves down slightly.
#   - The four "Ton/h" labels become "Ton/d" labels with new values.
#   - The "100 MW" label becomes a two-line "Potencia / Energía" label
#     (14pt) and grows taller.

$EMU_PER_POINT = 12700.0

# The COM layer stores Left/Top/Width/Height as single-precision floats, so a
# plain "$emu / $EMU_PER_POINT" tends to truncate down by ~1 EMU once it is
# converted back on save. Nudging by a fraction of an EMU before the divide
# keeps the round-trip exact.
function EmuToPt {
    param([double]$Emu)
    return ($Emu + 0.5) / $EMU_PER_POINT
}

function Get-ShapeById {
    param($Slide, [int]$Id)
    for ($i = 1; $i -le $Slide.Shapes.Count; $i++) {
        $candidate = $Slide.Shapes.Item($i)
        if ($candidate.Id -eq $Id) {
            return $candidate
        }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape id 10 "CuadroTexto 9" : NATURAL GAS ------------------------------
$shp = Get-ShapeById $s 10
$shp.TextFrame.TextRange.Text = "NATURAL GAS = 21.5 BFOE"
$shp.TextFrame.TextRange.Font.Size = 12
$shp.Left = EmuToPt 176549
$shp.Top = EmuToPt 1114463
$shp.Width = EmuToPt 1854418
$shp.Height = EmuToPt 276999

# --- Shape id 12 "CuadroTexto 11" : FUEL GAS --------------------------------
$shp = Get-ShapeById $s 12
$shp.TextFrame.TextRange.Text = "FUEL GAS = 504.3 BFOE"
$shp.TextFrame.TextRange.Font.Size = 12
$shp.Left = EmuToPt 176549
$shp.Top = EmuToPt 1590152
$shp.Width = EmuToPt 1653273
$shp.Height = EmuToPt 276999

# --- Shape id 14 connector between FUEL GAS and FLEXIGAS --------------------
$shp = Get-ShapeById $s 14
$shp.Top = EmuToPt 3033261

# --- Shape id 15 "CuadroTexto 14" : FLEXIGAS --------------------------------
$shp = Get-ShapeById $s 15
$shp.TextFrame.TextRange.Text = "FLEXIGAS = 3,133.8 BFOE"
$shp.TextFrame.TextRange.Font.Size = 12
$shp.Width = EmuToPt 1754263
$shp.Height = EmuToPt 276999

# --- Shape id 16 "CuadroTexto 15" : LIGHT NAPHTHA ---------------------------
$shp = Get-ShapeById $s 16
$shp.TextFrame.TextRange.Text = "LIGHT NAPHTHA = 0 BBL"
$shp.TextFrame.TextRange.Font.Size = 12
$shp.Top = EmuToPt 2659437
$shp.Width = EmuToPt 1710725
$shp.Height = EmuToPt 276999

# --- Shape id 2 "CuadroTexto 1" ---------------------------------------------
$shp = Get-ShapeById $s 2
$shp.TextFrame.TextRange.Text = "15,359 Ton/d"

# --- Shape id 3 "CuadroTexto 2" ---------------------------------------------
$shp = Get-ShapeById $s 3
$shp.TextFrame.TextRange.Text = "6,826 Ton/d"

# --- Shape id 7 "CuadroTexto 6" ---------------------------------------------
$shp = Get-ShapeById $s 7
$shp.TextFrame.TextRange.Text = "2,160 Ton/d"

# --- Shape id 26 "CuadroTexto 25" -------------------------------------------
$shp = Get-ShapeById $s 26
$shp.TextFrame.TextRange.Text = "5,374 Ton/d"

# --- Shape id 28 "CuadroTexto 27" : 100 MW ----------------------------------
$shp = Get-ShapeById $s 28
$shp.TextFrame.TextRange.Text = "Potencia = 100 MW" + [char]13 + "Energía = 1,354 BFOE/d"
$shp.TextFrame.TextRange.Font.Size = 14
$shp.Height = EmuToPt 523220
